$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "Project Engineer"
$ws.Range("B2").Value = "Jobot"
$ws.Range("C2").Value = "Dallas, TX"
$ws.Range("D2").Value = "Full-Time"
$ws.Range("E2").Value = "If you are a Residential / Project Engineer with experience in Airport and Aviation projects apply today! Why join us? Do you want to work with some of the nations best Clients AND enjoy time at home ..."
$ws.Range("F2").Value = "https://www.ziprecruiter.com/k/l/AAIQtvAWKQoy4PvqQA41MHcCmniElZZD1lN0oRqCUi-XetXtCzjRlIhw9yvhdq3BAdOKBRC1qL3ayVX2ksBvvbk2kVBN0AhplNgrFA48-k1sbe-_3-IHNoxwl7r8w7q9Bw9H-5cKRXS_n6Av8YqZOZYfMQCEj9QqeF2jsmkeaXI-edWoeEOUAw"
$ws.Range("A3").Value = "Help Desk Level 2 Engineer - Work From Home"
$ws.Range("B3").Value = "Ten4 Technology Group LLC"
$ws.Range("C3").Value = "Dallas, TX"
$ws.Range("D3").Value = "Full-Time"
$ws.Range("E3").Value = "We are searching for an experienced Managed Service Provider Level 2 engineer. The right candidate will possess superior customer service and support skills with prior experience working for a ..."
$ws.Range("F3").Value = "https://www.ziprecruiter.com/k/l/AAIn40gHlGO--ZFuED1IPvE0UGaxVfYRIeEBFePuFiyy4yPdJWM6pBV2O8w22EgZ1IwSzv8s03EEwENSJDhr9g-iSC-YevlLUbxH3d4ZiakfTND0SE2JRLYXOltAKiCSjmkgrV7_YQ9zEgBp4Ie7MI7gW4xn0P6JExL3FWWyS7aSfd9oQ_3ZoA"
$ws.Range("A4").Value = "Hotel Maintenance Engineer"
$ws.Range("B4").Value = "Crowne Plaza"
$ws.Range("C4").Value = "Dallas, TX"
$ws.Range("D4").Value = "Full-Time"
$ws.Range("E4").Value = "The CROWNE PLAZA DALLAS DOWNTOWN, a 300-unit upscale hotel, has an immediate opening for a professional and friendly Engineer in our maintenance department. Join a great group! Be part of the team ..."
$ws.Range("F4").Value = "https://www.ziprecruiter.com/k/l/AALF87QDEMNwO_JHyr5IO6iDig4O5dVjTyfcvUAfxRzMAg5r-RPGYXhT2BioW9ij29pLkE9pfamfO_GHob0DwXzG2TBVx02o-EOsgu4GH_8lCWposNdxBHYCva1efq62-1wCqB3mDVAq73FxqGEjojyzKkb7XQpaKM7I3KTCMX4WcqfNFxsc4g"
$ws.Range("A5").Value = "Audio Engineer"
$ws.Range("B5").Value = "Yo 214 Studios"
$ws.Range("C5").Value = "Dallas, TX"
$ws.Range("D5").Value = "Contractor"
$ws.Range("E5").Value = "Engineering in a professional, world-class studio environment working on all kinds of audio related projects, including music sessions for albums, radio/TV spots, podcasts, audiobooks, and audio-post ..."
$ws.Range("F5").Value = "https://www.ziprecruiter.com/k/l/AAI1VMESubnTPdlzDGUj9o7FOzmRuSYtsCvCUfSg9F6x-rQQEzzXmCkkG4_wlxIKJA-hyvY2aeUfYBGebZvDA4njvvgiUfm85_cyHw6YTPh8tFbP7Gb5bmAXNamxf1EQtbi51S1ysjGmiicX3zpftxF5vzMyJxo8mwf7Da7okOerl8SzCAkk1Q"
$ws.Range("A6").Value = "Support Services - TIER 2/3 Support & Operations Engineer"
$ws.Range("B6").Value = "ENEA"
$ws.Range("C6").Value = "Plano, TX"
$ws.Range("D6").Value = "Full-Time"
$ws.Range("E6").Value = "Preferrally college degree in Computer Science, Computer Engineer, Applied Math or related field. * Technical certification such as CCNA is highly desirable."
$ws.Range("F6").Value = "https://www.ziprecruiter.com/k/l/AALXEoa3M6K2elVy08KpWygpMdsyqv5zig2A9wPhifRsKfuUSAxcOkiKrOAPmMBQ5_FG8kamHM0DXBbkLlERdChvSAQLbgFb4gJtb3vc_eAAZR9IdLuuMVc3LuBB8XSqE7U-Nk0NxS6tnd1gvR1xW3ljnlSEGzJWbQztnBcOG_BVIA8okjEAHw"
$ws.Range("A7").Value = "Field Service Engineer- Dallas"
$ws.Range("B7").Value = "Volt"
$ws.Range("C7").Value = "Dallas, TX"
$ws.Range("D7").Value = "Full-Time"
$ws.Range("E7").Value = "Field Service Engineer Direct Hire Pay `$40+ hourly rate Great medical / dental benefits as well as 401K with match Dallas - Fort Worth, TX VOLT is working with a world leading manufacturer of ..."
$ws.Range("F7").Value = "https://www.ziprecruiter.com/k/l/AALn8eX5nr21vyb_6A6JIzqpvTjRBkubB4sSIMesadAxLqRfyqkj0Omi6dCXGY0VanTfYik3p9sQgjsW7YNJX4yFKDA7x7FppoYN-O6uarZfzk7A0XmdxSkj0njKlpjAa619zTSuGGDnhw2QRanBmViAYOC4XpxHF80OQ8ro6yFBekEXKWJzmg"
$ws.Range("A8").Value = "Telecom Network Virtualization Engineer"
$ws.Range("B8").Value = "Way Foward Consulting"
$ws.Range("C8").Value = "Plano, TX"
$ws.Range("D8").Value = "Full-Time"
$ws.Range("E8").Value = "Virtualization specifically NFV- Network Function Virtualization. (Minimum 1 year experience) RAN Virtualization & OAMP/OSS Engineer We are looking for a RAN Virtualization & OAMP Engineer to join ..."
$ws.Range("F8").Value = "https://www.ziprecruiter.com/k/l/AAKIhCQMnLwqXPcApOdeWOlZtiKY3rLNMgo9uLyKcDHfg5v2Ato03NAmnBXZkiBKuwXE0mkaPWzLi_mMKpufFGh3-TPyRDJ0f38Sn0EJrzId9Tvi-VhVf8K8KnB7GTzYRQJK5SSQppAb-4XKKQ6LygtbhoWUv675SScFU2br9fNHuMcdQSblRg"
$ws.Range("A9").Value = "Quality Engineer"
$ws.Range("B9").Value = "The Fountain Group"
$ws.Range("C9").Value = "Plano, TX"
$ws.Range("D9").Value = "Full-Time"
$ws.Range("E9").Value = "We are a national staffing firm and are currently seeking a -----Quality Engineer for a prominent client of ours. This position is located in -- Plano, TX. Details for the position are as follows"
$ws.Range("F9").Value = "https://www.ziprecruiter.com/k/l/AALdH2A9R1TTsDz4fOp1_CQ_0fSQsRAiVpkD4C2y0oU0bdRxAgPDfaRH0npBQ0oFLNvPqLH1pLHLid03l7aQ3TLrLWfmAAqHSZ2qE4fByLRxRjSLjKcpkzf1F7kVg0tltFvixNdalULspI9r0DP7XfSUleqPa5FQUpkEAfCJjdv9cDJ4zpc1Uw"
$ws.Range("A10").Value = "Call Processing Engineer--Telecom"
$ws.Range("B10").Value = "Way Foward Consulting"
$ws.Range("C10").Value = "Plano, TX"
$ws.Range("D10").Value = "Full-Time"
$ws.Range("E10").Value = "... Engineer to join our Wireless network team. You will be a key contributor as a subject matter expert in a multi-functional team of R&D, Systems engineering, and product validation to introduce ..."
$ws.Range("F10").Value = "https://www.ziprecruiter.com/k/l/AAJjgx17NCmhJtXg2BTZ1a2-KUgSAj-np9jBlnxw4DxIp6mZIcjuHtIlaw7Bk0hqhk9nE8YTW_JNv1-UuOtvveId1xajpF12MkeFveitiNlgO1oQDz2Zf10oyE_L8HJyeN7HGVViNY8FcvOdiqnSMlNhwxfvzNesqCJdyOedeiUoa5Ra9tuZzQ"
$ws.Range("A11").Value = "Lead Civil Designer"
$ws.Range("B11").Value = "Exceed Engineering, LLC"
$ws.Range("C11").Value = "Dallas, TX"
$ws.Range("D11").Value = "Full-Time"
$ws.Range("E11").Value = "Exceed Engineering is a multi-discipline professional engineering firm that provides engineering services to municipalities, state agencies and the private sector. We are a small but busy and growing ..."
$ws.Range("F11").Value = "https://www.ziprecruiter.com/k/l/AAKwgFfQQypfdo9oMfk3LGk0SANk-Bu4B-ukvPareXaMRwiLzCc7wj60KdKxB1naJdpaVZsQis4z5_7aF2VKBWEDnYtPC5S1K6UgqF8k2d-t1WfoA9vek1iIA0Y1xyY2yiEDySfQ3uyYH_eoDQrCpxKbhlcaw0Jvf8OdnMLf0JdC9scwzMJSODU"
$ws.Range("D12").Value = "Full-Time"
$ws.Range("D13").Value = "Full-Time"
$ws.Range("D14").Value = "Full-Time"
$ws.Range("D15").Value = "Full-Time"
$ws.Range("D16").Value = "Contractor"
